# Append the new run-log row (row 15) below the last existing row (row 14),
# matching the style used by the rest of the data rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 15
$prevRow = $newRow - 1

# Copy formatting (style) from the previous data row down into the new row
# before writing values, so the new row's cells inherit style index 3
# (center/center alignment) just like every other data row.
$srcRange = $ws.Range("A$prevRow" + ":H$prevRow")
$dstRange = $ws.Range("A$newRow" + ":H$newRow")
$srcRange.Copy($dstRange)

$ws.Cells.Item($newRow, 1).Value = "2025-08-15 06:51:35 UTC"
$ws.Cells.Item($newRow, 2).Value = "2025-08-15 12:21:35 IST"
$ws.Cells.Item($newRow, 3).Value = "SKIPPED"
$ws.Cells.Item($newRow, 4).Value = "No change in PDF. Skipping download & Excel update."
$ws.Cells.Item($newRow, 5).Value = "https://nalcoindia.com/wp-content/uploads/2025/08/INGOT-15-08-2025.pdf"
$ws.Cells.Item($newRow, 6).Value = ""
$ws.Cells.Item($newRow, 7).Value = 0
$ws.Cells.Item($newRow, 8).Value = ""
